$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Mes actual emp")
$ws2.Range("C2").Value = 160.05
$ws2.Select()
$ws2.Range("C5").Select()

$ws1 = $wb.Worksheets.Item("Mes actual cel")
$ws1.Range("C2").Value = 74.155
$ws1.Select()
$ws1.Range("E6").Select()
